$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("draftpicks")

$picks = @(
    @("balco", "Jesse Winker"),
    @("ottawa", "Jonathan Schoop"),
    @("deener", "Mitch Moreland"),
    @("virginia", "Michael Taylor"),
    @("dembums", "Matt Olson"),
    @("dembums", "Jordan Walden"),
    @("rippe", "Eric Young"),
    @("bears", "Kyle Schwarber"),
    @("virginia", "David Dejesus"),
    @("bellevegas", "Steven Moya"),
    @("isotopes", "Patrick Corbin"),
    @("jobu", "Nick Swisher"),
    @("ottawa", "Daniel Norris"),
    @("balco", "Michael Conforto"),
    @("balco", "Wilmer Flores"),
    @("ottawa", "Jose Peraza"),
    @("deener", "Jeurys Familia"),
    @("deener", "Ichiro Suzuki"),
    @("isotopes", "Josh Johnson"),
    @("bellevegas", "Luis Severino"),
    @("marmaduke", "Jameson Taillon"),
    @("pasadena", "Cody Asche"),
    @("virginia", "Ike Davis"),
)

$startRow = 474
for ($i = 0; $i -lt $picks.Count; $i++) {
    $r = $startRow + $i
    $team = $picks[$i][0]
    $player = $picks[$i][1]
    $ws.Cells.Item($r, 1).Value = $team
    $ws.Cells.Item($r, 2).Value = $player
    $ws.Cells.Item($r, 3).Value = 0
    $ws.Cells.Item($r, 4).Value = "B"
}

$ws.Range("C475:C496").Select()
$excel.ActiveWindow.ScrollRow = 460
$excel.ActiveWindow.ScrollColumn = 1

$ws.PageSetup.Orientation = 1
